# Trade #27 closed at 2026-02-16 21:26:48 - leadlag DOWN +0.000%
# Appends the new trade row (row 23) to the "leadlag" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$ws.Range("A23").Value = 27
# "2026-02-16" reads as a date to Excel's auto-detection; the leading
# apostrophe forces it to stay literal text, matching the other Date
# cells in this column (e.g. B2:B22).
$ws.Range("B23").Value = "'2026-02-16"
$ws.Range("C23").Value = "21:26:48"
$ws.Range("D23").Value = "leadlag"
$ws.Range("E23").Value = "DOWN"
$ws.Range("F23").Value = 68960.55
$ws.Range("H23").Value = "OPEN"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0.75
$ws.Range("L23").Value = "Coinbase leading with -0.110% move"
$ws.Range("N23").Value = 0
